$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    # Use Range.Text assignment (instead of Find.Execute's Replace argument)
    # so straight apostrophes in $new are not mangled by smart-quote
    # autocorrection. Loop to catch every occurrence of $old in the doc.
    $rng = $d.Content
    $rng.Start = 0
    $guard = 0
    while ($rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)) {
        $rng.Text = $new
        $rng.Collapse(0)
        $rng.End = $d.Content.End
        $guard = $guard + 1
        if ($guard -gt 50) { break }
    }
}

Replace-Text "Play Aurora Slot for Free - Review of Aurora Online Slot" "Play Aurora Slot Free - Exciting Gameplay and Stunning Visuals"
Replace-Text "2 wild symbols for increased win potential" "Simple but rewarding gameplay with 5 reels and 30 active paylines"
Replace-Text "Free spin bonus round with additional wilds" "Two wild symbols and a bonus round with the potential to fill the game grid with wilds"
Replace-Text "Beautiful graphics showcasing the northern lights" "Maximum win potential of 736 times your bet"
Replace-Text "Consistent payouts with 30 active paylines" "Attractive graphics and animations with a unique northern lights theme"
Replace-Text "Maximum win potential may not be very high" "Average RTP rate of 96.08%"
Replace-Text "Average RTP rate in the market" "Doesn't offer massive payouts"
Replace-Text "Discover Aurora, a beautiful online slot game with 5 reels and 30 paylines. Play for free and find out about its gameplay, graphics, and bonuses." "Read our review of Aurora slot and discover the rewarding gameplay and stunning visuals. Play for free!"
